$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (was ECs->ECs self row, becomes ECs->FAPs row)
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.302615
$ws.Range("N2").Value = 18.907845
$ws.Range("O2").Value = 0.6149105851929696
$ws.Range("P2").Value = 0.6149105851929696
$ws.Range("Q2").Value = 2.06132065667
$ws.Range("R2").Value = 18.55188591003
$ws.Range("S2").Value = 0.6149105851929696
$ws.Range("T2").Value = 0.6149105851929696

# Update row 3 (was ECs->FAPs row, becomes ECs->MuSCs row)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("M3").Value = 3.947029666666667
$ws.Range("N3").Value = 11.841089
$ws.Range("O3").Value = 0.3850894148070304
$ws.Range("P3").Value = 0.3850894148070304
$ws.Range("Q3").Value = 1.290907628720667
$ws.Range("R3").Value = 11.618168658486
$ws.Range("S3").Value = 0.3850894148070304
$ws.Range("T3").Value = 0.3850894148070304

# Remove old row 4 (was ECs->MuSCs row, now redundant after shifting data up)
$ws.Rows.Item(4).Delete()
